# Insert a new data row at row 196 (shifting existing rows 196-251 down to
# 197-252) in the single worksheet, then populate the new row with its
# values. This mirrors the diff: dimension grows from A1:R251 to A1:R252,
# and a brand-new weekly price observation is inserted right before the
# row that used to be 196.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 196..251 down to 197..252, creating a blank row 196.
$ws.Rows(196).Insert()

# Fill in the newly inserted row 196 with the new record's data.
$ws.Range("A196").Value = 9
$ws.Range("B196").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C196").Value = "Metropolitana"
$ws.Range("D196").Value = 44782
$ws.Range("E196").Value = 13
$ws.Range("F196").Value = 100112026
$ws.Range("G196").Value = "Haba"
$ws.Range("H196").Value = "Sin especificar"
$ws.Range("I196").Value = "Primera"
$ws.Range("J196").Value = 60
$ws.Range("K196").Value = 15000
$ws.Range("L196").Value = 15000
$ws.Range("M196").Value = 15000
$ws.Range("N196").Value = "$/saco 25 kilos"
$ws.Range("O196").Value = "Provincia del Elquí"
$ws.Range("P196").Value = 600
$ws.Range("Q196").Value = 25
$ws.Range("R196").Value = "Hortaliza"
